$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values for each row, preserving text representation
$priceUpdates = [ordered]@{
    2 = "243.86"
    4 = "5.244"
    6 = "6.467"
    7 = "3.231"
    8 = "0.8087"
    9 = "0.8914"
    10 = "0.1394"
    11 = "0.07092"
    12 = "0.03105"
    13 = "0.03049"
    14 = "0.09324"
    15 = "3.839"
    16 = "0.001560"
    17 = "0.04713"
    18 = "0.0006016"
    19 = "0.006171"
    20 = "0.001256"
    21 = "0.004065"
    22 = "0.00008704"
    24 = "2.159"
    25 = "0.3181"
    26 = "0.1322"
    28 = "0.0002331"
    40 = "0.03789"
    41 = "0.006283"
    42 = "0.1051"
    43 = "0.002530"
    44 = "0.007820"
    45 = "0.00005326"
    46 = "0.00000000750"
    47 = "0.5356"
    48 = "0.003430"
    49 = "0.00002101"
    50 = "0.0002001"
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    $origStyle = $cell.Style
    $cell.Value = "'" + $priceUpdates[$row]
    $cell.Style = $origStyle
}

# Update Volume(1h) label (column E) text for specific rows
$ws.Cells.Item(18, 5).Value = "17OneONEWorstin24h"
$ws.Cells.Item(47, 5).Value = "46CoinbaseStockTokenCOIN"
